$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2416691", "Lettuce - Romaine Hearts", "4", "36.35", "145.40"),
    @("8474538", "Spinach - Baby Fresh", "2", "18.15", "36.30"),
    @("0683696", "Tomato - Fresh Sliced", "12", "31.68", "380.16"),
    @("4966146", "Veggie Mix Power Blend", "1", "27.99", "27.99")
)

$startRow = 36
$endRow = $startRow + $newRows.Length - 1

$ws.Range("A$startRow`:E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
